$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: reassign owner/date note in column C from "Nam 15/02" to "Dũng 15/02"
$ws.Range("C9").Value = "Dũng 15/02"

# Row 17: B17 gets progress placeholder "?" (keeps existing percent number format style)
$ws.Range("B17").Value = "?"

# Insert 5 new blank rows before the old row 33 (Design Pattern Applied...) so
# that the existing "Design Pattern" rows shift from 33/34 down to 38/39, with
# only the first 3 of the new rows (30-32) ending up populated below.
$ws.Range("A30:A34").EntireRow.Insert()

# New Feature section header (row 30), styled like the other yellow section headers
$ws.Range("A30").Value = "Feature"
$ws.Range("A30").Interior.Color = $ws.Range("A3").Interior.Color

# New feature detail rows
$ws.Range("A31").Value = "Người chơi gây dame cho quái"
$ws.Range("A32").Value = "Quái gây dame cho người chơi"

# Update selection to match the post-edit state
$ws.Range("A33").Select()
